# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Range, [string]$Text)
    # Force literal text even when the value looks like a number
    # (e.g. '0.120' or '67.10') so trailing zeros / separators survive.
    $Range.Value = "'" + $Text
    $Range.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") "41.337.81"
$ws.Range("E2").Value = "  -1.16%  "

# Row 3
Set-TextCell $ws.Range("D3") "2.185.62"
$ws.Range("E3").Value = "  -1.56%  "

# Row 4
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
Set-TextCell $ws.Range("D5") "249.41"
$ws.Range("E5").Value = "  -1.43%  "

# Row 6
Set-TextCell $ws.Range("D6") "0.618"
$ws.Range("E6").Value = "  -2.14%  "

# Row 7
Set-TextCell $ws.Range("D7") "67.10"
$ws.Range("E7").Value = "  -4.75%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
Set-TextCell $ws.Range("D9") "0.612"
$ws.Range("E9").Value = "  +1.51%  "

# Row 10
Set-TextCell $ws.Range("D10") "38.40"
$ws.Range("E10").Value = "  -3.05%  "

# Row 11
Set-TextCell $ws.Range("D11") "59.50"
$ws.Range("E11").Value = "  +1.80%  "

# Row 12
Set-TextCell $ws.Range("D12") "0.0936"
$ws.Range("E12").Value = "  -3.42%  "

# Row 13
Set-TextCell $ws.Range("D13") "6.98"
$ws.Range("E13").Value = "  -4.70%  "

# Row 14
$ws.Range("E14").Value = "  -1.55%  "

# Row 15
Set-TextCell $ws.Range("D15") "2.514.58"
$ws.Range("E15").Value = "  -1.28%  "

# Row 16
Set-TextCell $ws.Range("D16") "14.47"
$ws.Range("E16").Value = "  -3.83%  "

# Row 17
Set-TextCell $ws.Range("D17") "0.854"
$ws.Range("E17").Value = "  -4.86%  "

# Row 18
Set-TextCell $ws.Range("D18") "2.188.12"
$ws.Range("E18").Value = "  -1.08%  "

# Row 19
Set-TextCell $ws.Range("D19") "41.299.50"
$ws.Range("E19").Value = "  -1.05%  "

# Row 20
Set-TextCell $ws.Range("D20") "0.0₃0950"
$ws.Range("E20").Value = "  -1.81%  "

# Row 21
Set-TextCell $ws.Range("D21") "71.87"
$ws.Range("E21").Value = "  -1.16%  "

# Row 22
$ws.Range("E22").Value = "  -3.11%  "

# Row 23
Set-TextCell $ws.Range("D23") "230.79"
$ws.Range("E23").Value = "  -1.93%  "

# Row 24
$ws.Range("E24").Value = "  -1.90%  "

# Row 25
$ws.Range("E25").Value = "  -5.36%  "

# Row 26
$ws.Range("E26").Value = "  +0.19%  "

# Row 27
Set-TextCell $ws.Range("D27") "11.25"
$ws.Range("E27").Value = "  -6.86%  "

# Row 28
$ws.Range("E28").Value = "  -5.31%  "

# Row 29
Set-TextCell $ws.Range("D29") "3.68"
$ws.Range("E29").Value = "  -2.94%  "

# Row 30
Set-TextCell $ws.Range("B30") "Monero"
Set-TextCell $ws.Range("C30") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws.Range("D30") "166.84"
$ws.Range("E30").Value = "  -2.91%  "

# Row 31
Set-TextCell $ws.Range("B31") "Toncoin"
Set-TextCell $ws.Range("C31") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell $ws.Range("D31") "2.07"
$ws.Range("E31").Value = "  -5.70%  "

# Row 32
Set-TextCell $ws.Range("D32") "20.22"
$ws.Range("E32").Value = "  -3.10%  "

# Row 33
Set-TextCell $ws.Range("D33") "0.0787"
$ws.Range("E33").Value = "  +5.21%  "

# Row 34
Set-TextCell $ws.Range("B34") "InternetComputer(DFINITY)"
Set-TextCell $ws.Range("C34") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell $ws.Range("D34") "5.84"
$ws.Range("E34").Value = "  +3.45%  "

# Row 35
Set-TextCell $ws.Range("B35") "Kaspa"
Set-TextCell $ws.Range("C35") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell $ws.Range("D35") "0.120"
$ws.Range("E35").Value = "  -3.05%  "

# Row 36
$ws.Range("E36").Value = "  -2.25%  "

# Row 37
$ws.Range("E37").Value = "  +3.93%  "

# Row 38
Set-TextCell $ws.Range("D38") "25.77"
$ws.Range("E38").Value = "  -0.30%  "

# Row 39
Set-TextCell $ws.Range("D39") "4.55"
$ws.Range("E39").Value = "  -3.39%  "

# Row 40
$ws.Range("E40").Value = "  -0.28%  "

# Row 41
Set-TextCell $ws.Range("D41") "2.20"
$ws.Range("E41").Value = "  -3.37%  "

# Row 42
$ws.Range("E42").Value = "  -6.12%  "

# Row 43
Set-TextCell $ws.Range("D43") "5.08"
$ws.Range("E43").Value = "  +5.67%  "

# Row 44
Set-TextCell $ws.Range("D44") "11.74"
$ws.Range("E44").Value = "  -5.25%  "

# Row 45
Set-TextCell $ws.Range("D45") "61.35"
$ws.Range("E45").Value = "  -5.93%  "

# Row 46
Set-TextCell $ws.Range("D46") "0.194"
$ws.Range("E46").Value = "  -5.61%  "

# Row 47
$ws.Range("E47").Value = "  -2.50%  "

# Row 48
Set-TextCell $ws.Range("D48") "8.48"
$ws.Range("E48").Value = "  -4.36%  "

# Row 49
$ws.Range("E49").Value = "  -0.24%  "

# Row 50
$ws.Range("E50").Value = "  -1.92%  "

# Row 51
$ws.Range("E51").Value = "  +5.40%  "
